$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New parameter values (lambdas) ---
# Row 4 = "investor growers"
$ws.Range("W4").Value = 2
$ws.Range("Z4").Value = 1

# Row 6 = "investor growers (white area)"
$ws.Range("Z6").Value = 1

# Row 26 = "Division of Water Quality (SWRCB)"
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 1

# Row 27 = "Groundwater Management (SWRCB)"
$ws.Range("C27").Value = 0.5
$ws.Range("D27").Value = 0.5
$ws.Range("E27").Value = 0.5
$ws.Range("F27").Value = 0.5

# Row 32 = "central valley water board"
$ws.Range("C32").Value = 0.5
$ws.Range("D32").Value = 0.5
$ws.Range("E32").Value = 0.5
$ws.Range("F32").Value = 0.5

# --- View / selection changes (comparing) ---
$ws.Application.ActiveWindow.SmallScroll(0, 0)
$ws.Range("H25").Select()
$excel.ActiveWindow.Panes.Item(1).ScrollRow = 1
$excel.ActiveWindow.Panes.Item(1).ScrollColumn = 1
$excel.ActiveWindow.Panes.Item(4).ScrollRow = 19
$excel.ActiveWindow.Panes.Item(4).ScrollColumn = 2
$excel.ActiveWindow.Panes.Item(4).Activate()
$ws.Range("H25").Select()
